$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.06187484444853
$ws.Range("D2").Value = 1.059904106548607
$ws.Range("E2").Value = 1.067592828746729
$ws.Range("F2").Value = 1.078726521730898
$ws.Range("I2").Value = 1.054818583205507
$ws.Range("J2").Value = 1.066848271596514
$ws.Range("K2").Value = 1.062632340361056
$ws.Range("L2").Value = 1.070300276584162
$ws.Range("M2").Value = 1.081404437849588
$ws.Range("N2").Value = 1.068363318511389
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.0631022890688
$ws.Range("D3").Value = 1.060845549766199
$ws.Range("E3").Value = 1.068718480113534
$ws.Range("F3").Value = 1.079998183666202
$ws.Range("I3").Value = 1.055257341296082
$ws.Range("J3").Value = 1.06772860528943
$ws.Range("K3").Value = 1.063387783049298
$ws.Range("L3").Value = 1.071240972104949
$ws.Range("M3").Value = 1.082492931420566
$ws.Range("N3").Value = 1.069244902379124
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.063896055791924
$ws.Range("D4").Value = 1.061454220259673
$ws.Range("E4").Value = 1.06944662483629
$ws.Range("F4").Value = 1.080821102506293
$ws.Range("I4").Value = 1.05553971672708
$ws.Range("J4").Value = 1.06829722493143
$ws.Range("K4").Value = 1.063875473487554
$ws.Range("L4").Value = 1.071848848568723
$ws.Range("M4").Value = 1.083196764576355
$ws.Range("N4").Value = 1.069814329526236
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.064229644730907
$ws.Range("D5").Value = 1.061709985319993
$ws.Range("E5").Value = 1.069752684222931
$ws.Range("F5").Value = 1.081167076020885
$ws.Range("I5").Value = 1.055658061852099
$ws.Range("J5").Value = 1.068536031138834
$ws.Range("K5").Value = 1.064080228903175
$ws.Range("L5").Value = 1.07210420595431
$ws.Range("M5").Value = 1.083492539341206
$ws.Range("N5").Value = 1.070053474865847
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.064285649366573
$ws.Range("D6").Value = 1.061752922395395
$ws.Range("E6").Value = 1.069804069901926
$ws.Range("F6").Value = 1.081225167636237
$ws.Range("I6").Value = 1.055677911111231
$ws.Range("J6").Value = 1.06857611363527
$ws.Range("K6").Value = 1.064114592462276
$ws.Range("L6").Value = 1.072147070259612
$ws.Range("M6").Value = 1.083542194435101
$ws.Range("N6").Value = 1.070093614284026
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.063900513655186
$ws.Range("D7").Value = 1.061457638275041
$ws.Range("E7").Value = 1.069450714621585
$ws.Range("F7").Value = 1.080825725346359
$ws.Range("I7").Value = 1.055541299495535
$ws.Range("J7").Value = 1.068300416818641
$ws.Range("K7").Value = 1.063878210497416
$ws.Range("L7").Value = 1.071852261426453
$ws.Range("M7").Value = 1.083200717189515
$ws.Range("N7").Value = 1.069817525946294
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.062289764292711
$ws.Range("D8").Value = 1.060222377216848
$ws.Range("E8").Value = 1.067973295909551
$ws.Range("F8").Value = 1.079156272866027
$ws.Range("I8").Value = 1.054967180960194
$ws.Range("J8").Value = 1.067145995601084
$ws.Range("K8").Value = 1.062887880784181
$ws.Range("L8").Value = 1.070618359039483
$ws.Range("M8").Value = 1.08177240292094
$ws.Range("N8").Value = 1.068661465318194
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.059447694758349
$ws.Range("D9").Value = 1.058041759248937
$ws.Range("E9").Value = 1.065368081159617
$ws.Range("F9").Value = 1.076214918126991
$ws.Range("I9").Value = 1.053943750467728
$ws.Range("J9").Value = 1.065103914335666
$ws.Range("K9").Value = 1.06113407162014
$ws.Range("L9").Value = 1.068437741477806
$ws.Range("M9").Value = 1.07925165877651
$ws.Range("N9").Value = 1.066616484063149
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.057550324916504
$ws.Range("D10").Value = 1.056585294785479
$ws.Range("E10").Value = 1.063629943194337
$ws.Range("F10").Value = 1.0742541769611
$ws.Range("I10").Value = 1.053253500233152
$ws.Range("J10").Value = 1.063737154782727
$ws.Range("K10").Value = 1.059958921599225
$ws.Range("L10").Value = 1.066979643335733
$ws.Range("M10").Value = 1.077568443909163
$ws.Range("N10").Value = 1.065247783554876
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.056728077850032
$ws.Range("D11").Value = 1.055953966794998
$ws.Range("E11").Value = 1.062876970395263
$ws.Range("F11").Value = 1.07340516320877
$ws.Range("I11").Value = 1.052952712524458
$ws.Range("J11").Value = 1.06314403606772
$ws.Range("K11").Value = 1.059448640432407
$ws.Range("L11").Value = 1.06634721505324
$ws.Range("M11").Value = 1.076838922026789
$ws.Range("N11").Value = 1.064653822543261
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.056422554413733
$ws.Range("D12").Value = 1.055719361284364
$ws.Range("E12").Value = 1.062597228313108
$ws.Range("F12").Value = 1.073089799070079
$ws.Range("I12").Value = 1.05284069921468
$ws.Range("J12").Value = 1.062923527599802
$ws.Range("K12").Value = 1.059258882251153
$ws.Range("L12").Value = 1.066112141361995
$ws.Range("M12").Value = 1.076567840775501
$ws.Range("N12").Value = 1.064433000928026
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.056488094995937
$ws.Range("D13").Value = 1.055769689613215
$ws.Range("E13").Value = 1.062657236439697
$ws.Range("F13").Value = 1.073157445909877
$ws.Range("I13").Value = 1.052864739472706
$ws.Range("J13").Value = 1.062970836396913
$ws.Range("K13").Value = 1.059299595873468
$ws.Range("L13").Value = 1.066162572816197
$ws.Range("M13").Value = 1.076625993383938
$ws.Range("N13").Value = 1.064480376909056
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.056702825332381
$ws.Range("D14").Value = 1.055934576324333
$ws.Range("E14").Value = 1.062853847951385
$ws.Range("F14").Value = 1.07337909517964
$ws.Range("I14").Value = 1.052943459338307
$ws.Range("J14").Value = 1.063125812815191
$ws.Range("K14").Value = 1.059432959399121
$ws.Range("L14").Value = 1.066327787098681
$ws.Range("M14").Value = 1.076816516508271
$ws.Range("N14").Value = 1.064635573411624
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.05683511387082
$ws.Range("D15").Value = 1.056036154879443
$ws.Range("E15").Value = 1.062974979509432
$ws.Range("F15").Value = 1.073515660182501
$ws.Range("I15").Value = 1.052991923131747
$ws.Range("J15").Value = 1.063221272634095
$ws.Range("K15").Value = 1.059515100259289
$ws.Range("L15").Value = 1.066429559585655
$ws.Range("M15").Value = 1.076933890195885
$ws.Range("N15").Value = 1.06473116879442
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.057604880294551
$ws.Range("D16").Value = 1.056627179704883
$ws.Range("E16").Value = 1.063679908005049
$ws.Range("F16").Value = 1.074310522984394
$ws.Range("I16").Value = 1.053273422289496
$ws.Range("J16").Value = 1.0637764904961
$ws.Range("K16").Value = 1.0599927569356
$ws.Range("L16").Value = 1.067021592961133
$ws.Range("M16").Value = 1.077616845371548
$ws.Range("N16").Value = 1.065287175129473
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.058087551963474
$ws.Range("D17").Value = 1.056997733628832
$ws.Range("E17").Value = 1.064121996388964
$ws.Range("F17").Value = 1.074809117583887
$ws.Range("I17").Value = 1.053449488466894
$ws.Range("J17").Value = 1.06412441375719
$ws.Range("K17").Value = 1.06029199339149
$ws.Range("L17").Value = 1.067392674010327
$ws.Range("M17").Value = 1.07804506172114
$ws.Range("N17").Value = 1.065635592481502
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.058369021418558
$ws.Range("D18").Value = 1.057213807000794
$ws.Range("E18").Value = 1.064379825678941
$ws.Range("F18").Value = 1.075099939569239
$ws.Range("I18").Value = 1.053552001184509
$ws.Range("J18").Value = 1.064327225945589
$ws.Range("K18").Value = 1.060466394760449
$ws.Range("L18").Value = 1.067609016967062
$ws.Range("M18").Value = 1.078294767746184
$ws.Range("N18").Value = 1.065838692686471
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.058464984336971
$ws.Range("D19").Value = 1.057287471551757
$ws.Range("E19").Value = 1.064467733208385
$ws.Range("F19").Value = 1.075199102462565
$ws.Range("I19").Value = 1.053586924248152
$ws.Range("J19").Value = 1.064396358424085
$ws.Range("K19").Value = 1.060525837730636
$ws.Range("L19").Value = 1.067682767028424
$ws.Range("M19").Value = 1.078379900038326
$ws.Range("N19").Value = 1.065907923341016
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.05803577258377
$ws.Range("D20").Value = 1.056957983380107
$ws.Range("E20").Value = 1.064074567959683
$ws.Range("F20").Value = 1.074755623102615
$ws.Range("I20").Value = 1.053430617242959
$ws.Range("J20").Value = 1.064087097887288
$ws.Range("K20").Value = 1.060259902458638
$ws.Range("L20").Value = 1.067352871127355
$ws.Range("M20").Value = 1.077999124934736
$ws.Range("N20").Value = 1.065598223618784
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.05663959550343
$ws.Range("D21").Value = 1.055886024131581
$ws.Range("E21").Value = 1.062795952282926
$ws.Range("F21").Value = 1.073313825088007
$ws.Range("I21").Value = 1.05292028624208
$ws.Range("J21").Value = 1.063080181598457
$ws.Range("K21").Value = 1.059393693169544
$ws.Range("L21").Value = 1.066279140089996
$ws.Range("M21").Value = 1.076760415096238
$ws.Range("N21").Value = 1.064589877393328
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.055761158774517
$ws.Range("D22").Value = 1.055211448779578
$ws.Range("E22").Value = 1.061991719666336
$ws.Range("F22").Value = 1.072407292927114
$ws.Range("I22").Value = 1.052597757707885
$ws.Range("J22").Value = 1.062445948122734
$ws.Range("K22").Value = 1.058847816474904
$ws.Range("L22").Value = 1.065603106368919
$ws.Range("M22").Value = 1.07598098418666
$ws.Range("N22").Value = 1.063954743233322
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.056226893067246
$ws.Range("D23").Value = 1.055569110602579
$ws.Range("E23").Value = 1.062418089313431
$ws.Range("F23").Value = 1.072887865229707
$ws.Range("I23").Value = 1.052768894233802
$ws.Range("J23").Value = 1.062782276525689
$ws.Range("K23").Value = 1.059137315705211
$ws.Range("L23").Value = 1.065961574095138
$ws.Range("M23").Value = 1.076394233387079
$ws.Range("N23").Value = 1.064291549261186
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.058059169656702
$ws.Range("D24").Value = 1.056975945005711
$ws.Range("E24").Value = 1.064095998928522
$ws.Range("F24").Value = 1.074779794955479
$ws.Range("I24").Value = 1.053439144904371
$ws.Range("J24").Value = 1.064103959710624
$ws.Range("K24").Value = 1.060274403396864
$ws.Range("L24").Value = 1.067370856653098
$ws.Range("M24").Value = 1.078019881993041
$ws.Range("N24").Value = 1.065615109387843
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.06018289487085
$ws.Range("D25").Value = 1.058605974883285
$ws.Range("E25").Value = 1.066041817749259
$ws.Range("F25").Value = 1.076975290752448
$ws.Range("I25").Value = 1.054209731221827
$ws.Range("J25").Value = 1.065632780914821
$ws.Range("K25").Value = 1.061588515105306
$ws.Range("L25").Value = 1.069002243211526
$ws.Range("M25").Value = 1.079903802276465
$ws.Range("N25").Value = 1.067146101693509
